$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values (column B) and units (column C) for the frame specification rows.
# Column A labels stay the same text except row 5, which is renamed.
$ws.Range("B1").Value = 400
$ws.Range("C1").Value = "mm"

$ws.Range("B2").Value = 150
$ws.Range("C2").Value = "mm"

$ws.Range("B3").Value = 100
$ws.Range("C3").Value = "mm"

$ws.Range("B4").Value = 150
$ws.Range("C4").Value = "mm"

$ws.Range("A5").Value = "suppourtPanelThickness"
$ws.Range("B5").Value = 0.09375
$ws.Range("C5").Value = "in"

# Move the active selection to B2.
$ws.Range("B2").Select() | Out-Null

# Set page setup to A4 / portrait orientation.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1 | Out-Null
